$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh: updated Price (D) and Volume 1h (E) columns,
# plus a position swap of BabyDogeCoin/OKB (rows 46-47) with fresh values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.085.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.27%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.741.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.72%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.56%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  +1.83%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.740.47"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.71%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.148"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.03%  "

$ws.Range("E11").Value = "  +5.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.37"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.66%  "

$ws.Range("E13").Value = "  -0.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.83%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.239.03"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.56%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000193"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.33%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.069.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.44%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.722.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.78%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "378.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.48%  "

$ws.Range("E21").Value = "  +4.92%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.55%  "

$ws.Range("E24").Value = "  +3.47%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.68%  "

$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.16%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000107"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "591.63"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.78%  "

$ws.Range("E31").Value = "  -0.11%  "

$ws.Range("E32").Value = "  +5.56%  "

$ws.Range("E33").Value = "  +5.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.01"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.23%  "

$ws.Range("E35").Value = "  +4.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.99%  "

$ws.Range("E37").Value = "  +0.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "162.64"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.15%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "20.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.12%  "

$ws.Range("E40").Value = "  +4.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.74%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.43%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.17%  "

$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0313"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.48%  "

$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "41.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.78%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.610"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.47%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "156.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.55%  "

$ws.Range("E50").Value = "  +4.87%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.07%  "
